# RPA datasets push 2024-05-04
# Applies the dataset refresh to both sheets of the IB strategy workbook:
#  - Sheet "01_IB전략컨설팅부": remove the SK증권제11호스팩 entry (old row 15)
#  - Sheet "02_38커뮤니케이션(최근일자기준)": refresh the rolling 수요예측 list
#    (drop the three oldest entries, add three new ones, and patch a few
#    values that were finalised since the last pull)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Sheet 1 (01_IB전략컨설팅부): drop the SK증권제11호스팩 row entirely —
# everything below it (the BNK / 비엔케이제2호스팩 row) shifts up.
# ---------------------------------------------------------------------
$ws1.Rows.Item(15).Delete()

# ---------------------------------------------------------------------
# Sheet 2 (02_38커뮤니케이션(최근일자기준)): roll the demand-forecast table
# forward. Work from the bottom of the sheet upward so earlier row
# numbers used below stay valid.
# ---------------------------------------------------------------------

# Drop the three oldest entries that have rolled off the bottom of the list.
$ws2.Range("A19:F21").EntireRow.Delete()

# 아이씨티케이 and KB스팩28호 now have confirmed offering prices.
$ws2.Cells.Item(13, 4).Value = "20000"
$ws2.Cells.Item(12, 4).Value = "2000"

# Two new spac/IPO entries land right after 미래에셋비전스팩5호.
$ws2.Rows.Item(7).Insert()
$ws2.Cells.Item(7, 1).Value = "그리드위즈"
$ws2.Cells.Item(7, 2).Value = "2024.05.23~05.29"
$ws2.Cells.Item(7, 3).Value = "34,000~40,000"
$ws2.Cells.Item(7, 4).Value = "-"
$ws2.Cells.Item(7, 5).Value = "47600"
$ws2.Cells.Item(7, 6).Value = "삼성증권"

$ws2.Rows.Item(8).Insert()
$ws2.Cells.Item(8, 1).Value = "이노스페이스"
$ws2.Cells.Item(8, 2).Value = "2024.05.23~05.29"
$ws2.Cells.Item(8, 3).Value = "36,400~45,600"
$ws2.Cells.Item(8, 4).Value = "-"
$ws2.Cells.Item(8, 5).Value = "48412"
$ws2.Cells.Item(8, 6).Value = "미래에셋증권,신한투자증권"

# 씨어스테크놀로지's demand-forecast window moved.
$ws2.Cells.Item(6, 2).Value = "2024.05.27~05.31"

# New entry 한국스팩14호 lands right after 하이젠알앤엠.
$ws2.Rows.Item(3).Insert()
$ws2.Cells.Item(3, 1).Value = "한국스팩14호"
$ws2.Cells.Item(3, 2).Value = "2024.06.03~06.04"
$ws2.Cells.Item(3, 3).Value = "2,000~2,000"
$ws2.Cells.Item(3, 4).Value = "-"
$ws2.Cells.Item(3, 5).Value = "8000"
$ws2.Cells.Item(3, 6).Value = "한국투자증권"
